# TC04_Canine_Filter_Diagnosis-Melanoma.xlsx — "10 icdc scripts for jenkins"
#
# The "startup" sheet's FilesTab row (B4) holds the Cypher query used to
# build the Files report. Drop the `File Type` and `Breed` columns from its
# RETURN clause (they are not part of this report).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$fileQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.disease_term IN ['Melanoma']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value2 = $fileQuery

# Author last had row 4 in view with B4 selected when the workbook was saved.
[void]$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
